$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "Bug in Current Balance..." status Pending -> Done ---
$ws.Range("I5:K5").Copy()
$ws.Range("I13:K13").PasteSpecial(-4122)
$ws.Range("I13:K13").Value = "Done"
$ws.Rows(13).RowHeight = 16

# --- Row 15: "Get entries by month filter" status Pending -> Done ---
$ws.Range("I5:K5").Copy()
$ws.Range("I15:K15").PasteSpecial(-4122)
$ws.Range("I15:K15").Value = "Done"
$ws.Rows(15).RowHeight = 16

# --- Row 16: new backlog entry ---
$ws.Range("A15:D15").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)
$ws.Range("A16:D16").Value = "Get All entries for a product and sort by entry date"

$ws.Range("I5:K5").Copy()
$ws.Range("I16:K16").PasteSpecial(-4122)
$ws.Range("I16:K16").Value = "Done"
$ws.Rows(16).RowHeight = 49

$excel.CutCopyMode = 0

# Restore the active selection to match the authored state
$ws.Range("M13").Select()
